# Generate Report for Handoff
# Swap the two tracked files' row positions (66fa49ac.md <-> f47fb380.md) on every
# sheet, and record that 66fa49ac.md has now moved from "In Translation" to
# "Ready for handoff" with a fresh handoff file/date, while f47fb380.md stays
# "In Translation" (unchanged details).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "f47fb380-86a0-4d9d-a260-925a5d454aaa.md"
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("D2").Value = "2016-13-12 10:13:23"

$wsOverview.Range("A3").Value = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-14-12 10:14:13"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "f47fb380-86a0-4d9d-a260-925a5d454aaa.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | ... | Latest Handback DateTime |
# Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "f47fb380-86a0-4d9d-a260-925a5d454aaa.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("D2").Value = "f47fb380-86a0-4d9d-a260-925a5d454aaa.029e7c78ae154090c69153c646152f02bae5b7ee.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-12 10:12:57"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.42ae7acb600d536657d570ae13f29341b33eafa2.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-12 10:14:10"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "f47fb380-86a0-4d9d-a260-925a5d454aaa.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "f47fb380-86a0-4d9d-a260-925a5d454aaa.029e7c78ae154090c69153c646152f02bae5b7ee.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.42ae7acb600d536657d570ae13f29341b33eafa2.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de": same column layout as "zh-cn"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "f47fb380-86a0-4d9d-a260-925a5d454aaa.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("D2").Value = "f47fb380-86a0-4d9d-a260-925a5d454aaa.029e7c78ae154090c69153c646152f02bae5b7ee.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-12 10:13:23"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.42ae7acb600d536657d570ae13f29341b33eafa2.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-12 10:14:13"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "f47fb380-86a0-4d9d-a260-925a5d454aaa.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "f47fb380-86a0-4d9d-a260-925a5d454aaa.029e7c78ae154090c69153c646152f02bae5b7ee.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "66fa49ac-8554-4129-8ee0-b9fd8dfd5cc5.42ae7acb600d536657d570ae13f29341b33eafa2.de-de.xlf"
    }
}
